# "save result to DB": a second quiz submission came in (Vu Thi Huyen),
# so we record her row on the Summary sheet, refresh the submission
# timestamp that was already there for Yen Tuan Phong, and append her
# full per-question answer sheet (same shape as the existing per-student
# sheets) as a new worksheet tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Summary sheet: bump Yen Tuan Phong's CompletedTime and add the
#    new row for Vu Thi Huyen.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = "05-03-2023 00:23:36"

$summary.Range("A3").Value = "Vũ Thị Huyền"
# StudentScore is stored as text in this workbook ("0", "1", ...), not a
# number -- the leading apostrophe forces Excel to keep it as text
# instead of silently converting it to a numeric 1.
$summary.Range("B3").Value = "'1"
$summary.Range("C3").Value = "05-03-2023 00:27:57"

# ---------------------------------------------------------------
# 2) New worksheet for Vu Thi Huyen's answers, appended after the last
#    existing tab (same layout used for Yen Tuan Phong's sheet).
# ---------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Vũ Thị Huyền"

# Column widths (A..H) = 25,25,25,25,8,25,14,11 character-widths.
# Excel's ColumnWidth property adds ~0.8333 to whatever raw "width" ends
# up in the XML, so subtract that back off to land on the exact target.
$colWidths = @(25, 25, 25, 25, 8, 25, 14, 11)
for ($i = 0; $i -lt $colWidths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $colWidths[$i] - 0.8333333333333334
}

# Header row
$headers = @("Answer_A", "Answer_B", "Answer_C", "Answer_D", "Num", "Question", "Answer_Correct", "Answer_User")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$headerRange = $ws.Range("A1:H1")
$headerRange.Font.Name = "arial"
$headerRange.Font.Size = 10
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 0xC7A52E
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.Borders.Color = 0x000000

# Question rows
$ws.Range("A2").Value = "Hyper Text Preprocessor A"
$ws.Range("B2").Value = "Hyper Text Preprocessor B"
$ws.Range("C2").Value = "Hyper Text Preprocessor C"
$ws.Range("D2").Value = "Hyper Text Preprocessor D"
$ws.Range("E2").Value = "'1"
$ws.Range("F2").Value = "What does HTML stand for?"
$ws.Range("G2").Value = "Answer_A"
$ws.Range("H2").Value = "Answer_A"

$ws.Range("A3").Value = "Hyper Text A"
$ws.Range("C3").Value = "Hyper Text C"
$ws.Range("E3").Value = "'2"
$ws.Range("F3").Value = "What does CSS stand for?"
$ws.Range("G3").Value = "Answer_C"

$ws.Range("A4").Value = "JavaScript A"
$ws.Range("B4").Value = "JavaScript B"
$ws.Range("C4").Value = "JavaScript C"
$ws.Range("D4").Value = "JavaScript D"
$ws.Range("E4").Value = "'3"
$ws.Range("F4").Value = "What does JS stand for?"
$ws.Range("G4").Value = "Answer_A"

# Restore the originally active sheet as the selected tab.
$summary.Activate()
